$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2").Value = "#! END_ROW true"
$ws.Range("J3").Value = "#! END_ROW"
